$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 613.93024
$ws.Range("I6").Value = 100
$ws.Range("J6").Value = 2555.4443
$ws.Range("K6").Value = 300
$ws.Range("L6").Value = 7666.3329
$ws.Range("M6").Value = -188
$ws.Range("N6").Value = -7890.3329
$ws.Range("H8").Value = 275
$ws.Range("I8").Value = 275
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 825
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -686
$ws.Range("N8").ClearContents()
$ws.Range("H62").Value = 2396.7856
$ws.Range("I62").Value = 1617.2222
$ws.Range("K62").Value = 1617.2222
$ws.Range("M62").Value = -993.2221999999999
$ws.Range("H65").Value = 2396.7856
$ws.Range("I65").Value = 1617.2222
$ws.Range("K65").Value = 8086.111
$ws.Range("M65").Value = -4966.111
$ws.Range("H107").Value = 1105.4546
$ws.Range("I107").Value = 1151.0588
$ws.Range("J107").Value = 950.4
$ws.Range("K107").Value = 1151.0588
$ws.Range("L107").Value = 950.4
$ws.Range("M107").Value = 768.9412
$ws.Range("N107").Value = -4790.4
$ws.Range("H112").Value = 5927.4526
$ws.Range("J112").Value = 6303.163
$ws.Range("L112").Value = 18909.489
$ws.Range("N112").Value = -21125.489
$ws.Range("H138").Value = 180188.3
$ws.Range("I138").Value = 2200.95
$ws.Range("J138").Value = 267011.38
$ws.Range("K138").Value = 6602.849999999999
$ws.Range("L138").Value = 801034.14
$ws.Range("M138").Value = -1462.849999999999
$ws.Range("N138").Value = -811314.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 21500
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H32").Value = 500672.56
$ws.Range("I32").Value = 640338.4399999999
$ws.Range("J32").Value = 11842
$ws.Range("K32").Value = 640338.4399999999
$ws.Range("L32").Value = 11842
$ws.Range("M32").Value = -640051.4399999999
$ws.Range("N32").Value = -12416
$ws.Range("H74").Value = 1869.0714
$ws.Range("I74").Value = 1705.1538
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 1705.1538
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -831.1538
$ws.Range("N74").Value = -5748
$ws.Range("H77").Value = 1869.0714
$ws.Range("I77").Value = 1705.1538
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 8525.769
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -4157.769
$ws.Range("N77").Value = -28736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 2675.1538
$ws.Range("I22").Value = 2766.16
$ws.Range("K22").Value = 2766.16
$ws.Range("M22").Value = -2593.16
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H134").Value = 2628.111
$ws.Range("I134").Value = 2015.2354
$ws.Range("J134").Value = 3670
$ws.Range("K134").Value = 6045.706200000001
$ws.Range("L134").Value = 11010
$ws.Range("M134").Value = -3510.706200000001
$ws.Range("N134").Value = -16080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.727272
$ws.Range("I7").Value = 43.57143
$ws.Range("J7").Value = 90.75
$ws.Range("K7").Value = 43.57143
$ws.Range("L7").Value = 90.75
$ws.Range("M7").Value = 69.42857000000001
$ws.Range("N7").Value = -316.75
$ws.Range("H31").Value = 1538.5
$ws.Range("I31").Value = 936.8333
$ws.Range("J31").Value = 3042.6667
$ws.Range("K31").Value = 936.8333
$ws.Range("L31").Value = 3042.6667
$ws.Range("M31").Value = -641.8333
$ws.Range("N31").Value = -3632.6667
$ws.Range("H34").Value = 1538.5
$ws.Range("I34").Value = 936.8333
$ws.Range("J34").Value = 3042.6667
$ws.Range("K34").Value = 936.8333
$ws.Range("L34").Value = 3042.6667
$ws.Range("M34").Value = -734.8333
$ws.Range("N34").Value = -3446.6667
$ws.Range("H43").Value = 42487.5
$ws.Range("J43").Value = 42487.5
$ws.Range("L43").Value = 42487.5
$ws.Range("N43").Value = -42855.5
$ws.Range("H58").Value = 1322.9474
$ws.Range("I58").Value = 852.4
$ws.Range("K58").Value = 852.4
$ws.Range("M58").Value = -649.4
$ws.Range("H94").Value = 1173.6
$ws.Range("I94").Value = 400
$ws.Range("J94").Value = 1367
$ws.Range("K94").Value = 400
$ws.Range("L94").Value = 1367
$ws.Range("M94").Value = 51
$ws.Range("N94").Value = -2269
$ws.Range("H101").Value = 42487.5
$ws.Range("J101").Value = 42487.5
$ws.Range("L101").Value = 42487.5
$ws.Range("N101").Value = -48977.5
$ws.Range("H132").Value = 10418762
$ws.Range("I132").Value = 1100.2858
$ws.Range("K132").Value = 3300.8574
$ws.Range("M132").Value = -770.8574000000003
$ws.Range("H136").Value = 1322.9474
$ws.Range("I136").Value = 852.4
$ws.Range("K136").Value = 2557.2
$ws.Range("M136").Value = -7.199999999999818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 109.375
$ws.Range("I2").Value = 25.714285
$ws.Range("J2").Value = 174.44444
$ws.Range("K2").Value = 154.28571
$ws.Range("L2").Value = 1046.66664
$ws.Range("M2").Value = -41.28570999999999
$ws.Range("N2").Value = -1272.66664
$ws.Range("H17").Value = 520
$ws.Range("I17").Value = 562.5
$ws.Range("J17").Value = 180
$ws.Range("K17").Value = 1687.5
$ws.Range("L17").Value = 540
$ws.Range("M17").Value = -1518.5
$ws.Range("N17").Value = -878
$ws.Range("H33").Value = 200249.8
$ws.Range("I33").Value = 33699.668
$ws.Range("J33").Value = 450075
$ws.Range("K33").Value = 202198.008
$ws.Range("L33").Value = 2700450
$ws.Range("M33").Value = -201915.008
$ws.Range("N33").Value = -2701016
$ws.Range("H113").Value = 2764.6667
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 3697
$ws.Range("K113").Value = 2700
$ws.Range("L113").Value = 11091
$ws.Range("M113").Value = -530
$ws.Range("N113").Value = -15431
$ws.Range("H121").Value = 1121.4255
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1121.4255
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 3364.2765
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5984.2765
$ws.Range("H131").Value = 962.34
$ws.Range("I131").Value = 297.5
$ws.Range("J131").Value = 1020.15216
$ws.Range("K131").Value = 892.5
$ws.Range("L131").Value = 3060.45648
$ws.Range("M131").Value = 4147.5
$ws.Range("N131").Value = -13140.45648

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 5000
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224
$ws.Range("H93").Value = 49998.46
$ws.Range("J93").Value = 49998.46
$ws.Range("L93").Value = 49998.46
$ws.Range("N93").Value = -53742.46
$ws.Range("H104").Value = 33167.75
$ws.Range("J104").Value = 33167.75
$ws.Range("L104").Value = 33167.75
$ws.Range("N104").Value = -40155.75
$ws.Range("H132").Value = 2495.4443
$ws.Range("I132").Value = 1730.6
$ws.Range("J132").Value = 3451.5
$ws.Range("K132").Value = 5191.799999999999
$ws.Range("L132").Value = 10354.5
$ws.Range("M132").Value = -2661.799999999999
$ws.Range("N132").Value = -15414.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 25333.666
$ws.Range("J2").Value = 39857.715
$ws.Range("L2").Value = 39857.715
$ws.Range("N2").Value = -40081.715
$ws.Range("H109").Value = 30000
$ws.Range("J109").Value = 30000
$ws.Range("L109").Value = 30000
$ws.Range("N109").Value = -32774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 77950
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 102633.336
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 102633.336
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -103881.336
$ws.Range("H65").Value = 77950
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 102633.336
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 513166.68
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -519406.68
$ws.Range("H105").Value = 100000
$ws.Range("J105").Value = 100000
$ws.Range("L105").Value = 100000
$ws.Range("N105").Value = -106988
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
